$d = $word.ActiveDocument

# Change 1: remove ", filtrando dinamicamente gli indirizzi IP" before the final period
$d.Content.Find.Execute(
    "IL Firewall controlla così il traffico in entrata ed in uscita, filtrando dinamicamente gli indirizzi IP.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IL Firewall controlla così il traffico in entrata ed in uscita.",
    2
)

# Change 2: append new sentences about WAF after the SMTP sentence
$d.Content.Find.Execute(
    "Questa zona ci serve per poter far raggiungere il server di servizi web (HTTP) e il server di servizi mail (SMTP).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questa zona ci serve per poter far raggiungere il server di servizi web (HTTP) e il server di servizi mail (SMTP). La DMZ è protetta dal WAF (Web Application Firewall) che filtra i contenuti dei pacchetti che riceve, attraverso il confronto della sua tabella con una tabella di terzi parti. Per esempio confronta il contenuto del pacchetto con i server di OWASP per capire se ci sono malware nel pacchetto. ",
    2
)

# Change 3: remove leading space and insert "anche" before "per l'IPS"
$d.Content.Find.Execute(
    " Prima di arrivare ai due server, i dati passano per l’IPS (Intrusion protection system).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prima di arrivare ai due server, i dati passano anche per l’IPS (Intrusion protection system).",
    2
)
